$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update angle (A) values for rows 26-37: 20 -> 17
$ws.Range("A26:A37").Value = 17

# Update starting angle offset B26: 18 -> 0
# (B27:B37 are formulas referencing the previous row, so they recalc automatically)
$ws.Range("B26").Value = 0

# Update the selected cell / view to B27
$ws.Range("B27").Select()
